$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(8)
$tr = $sh.TextFrame.TextRange

function Add-Run($para, $runText, $size) {
    $r = $para.InsertAfter($runText)
    $r.Font.Size = $size
    $r.LanguageID = "en-US"
    $r.Font.NameFarEast = "+mn-lt"
    $r.Font.NameComplexScript = "+mn-lt"
    [void]$r
}

# Paragraph 1
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = ""
[void](Add-Run $para1 "Genetski" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "algoritam" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "predstavlja" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "populacionu" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "metaheuristiku" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "inspirisanu" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "prirodnom" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "selekcijom" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "i" 2000)
[void](Add-Run $para1 " " 2000)
[void](Add-Run $para1 "evolucijom" 2000)
[void](Add-Run $para1 "." 2000)

# Paragraph 2
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = ""
[void](Add-Run $para2 "Svaki" 2000)
[void](Add-Run $para2 " " 2000)
[void](Add-Run $para2 "kandidat" 2000)
[void](Add-Run $para2 " za " 2000)
[void](Add-Run $para2 "rešenje" 2000)
[void](Add-Run $para2 " Nearest String " 2000)
[void](Add-Run $para2 "problema" 2000)
[void](Add-Run $para2 " " 2000)
[void](Add-Run $para2 "predstavlja" 2000)
[void](Add-Run $para2 " " 2000)
[void](Add-Run $para2 "jednu" 2000)
[void](Add-Run $para2 " " 2000)
[void](Add-Run $para2 "jedinku" 2000)
[void](Add-Run $para2 " u " 2000)
[void](Add-Run $para2 "populaciji" 2000)
[void](Add-Run $para2 "." 2000)

# Paragraph 3
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = ""
[void](Add-Run $para3 "Algoritam" 2000)
[void](Add-Run $para3 " se " 2000)
[void](Add-Run $para3 "sastoji" 2000)
[void](Add-Run $para3 " od " 2000)
[void](Add-Run $para3 "sledećih" 2000)
[void](Add-Run $para3 " " 2000)
[void](Add-Run $para3 "faza" 2000)
[void](Add-Run $para3 " (" 2000)
[void](Add-Run $para3 "koje" 2000)
[void](Add-Run $para3 " se " 2000)
[void](Add-Run $para3 "ponavljaju" 2000)
[void](Add-Run $para3 " " 2000)
[void](Add-Run $para3 "kroz" 2000)
[void](Add-Run $para3 " " 2000)
[void](Add-Run $para3 "unapred" 2000)
[void](Add-Run $para3 " " 2000)
[void](Add-Run $para3 "definisan" 2000)
[void](Add-Run $para3 " " 2000)
[void](Add-Run $para3 "broj" 2000)
[void](Add-Run $para3 " " 2000)
[void](Add-Run $para3 "generacija" 2000)
[void](Add-Run $para3 "):" 2000)

# Paragraph 4
$para4 = $tr.Paragraphs(4, 1)
$para4.Text = ""
[void](Add-Run $para4 "Inicijalizacija" 1600)
[void](Add-Run $para4 " " 1600)
[void](Add-Run $para4 "populacije" 1600)
[void](Add-Run $para4 "," 1600)

# Paragraph 5
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = ""
[void](Add-Run $para5 "Procena" 1600)
[void](Add-Run $para5 " " 1600)
[void](Add-Run $para5 "prilagođenosti" 1600)
[void](Add-Run $para5 " (fitness)," 1600)

# Paragraph 6
$para6 = $tr.Paragraphs(6, 1)
$para6.Text = ""
[void](Add-Run $para6 "Selekcija" 1600)
[void](Add-Run $para6 " " 1600)
[void](Add-Run $para6 "roditelja" 1600)
[void](Add-Run $para6 "," 1600)

# Paragraph 7
$para7 = $tr.Paragraphs(7, 1)
$para7.Text = ""
[void](Add-Run $para7 "Ukrštanje" 1600)
[void](Add-Run $para7 "," 1600)

# Paragraph 8
$para8 = $tr.Paragraphs(8, 1)
$para8.Text = ""
[void](Add-Run $para8 "Mutacija" 1600)
[void](Add-Run $para8 "," 1600)

# Paragraph 9
$para9 = $tr.Paragraphs(9, 1)
$para9.Text = ""
[void](Add-Run $para9 "Elitizam" 1600)
[void](Add-Run $para9 " " 1600)
[void](Add-Run $para9 "i" 1600)
[void](Add-Run $para9 " " 1600)
[void](Add-Run $para9 "formiranje" 1600)
[void](Add-Run $para9 " " 1600)
[void](Add-Run $para9 "nove" 1600)
[void](Add-Run $para9 " " 1600)
[void](Add-Run $para9 "generacije" 1600)
[void](Add-Run $para9 "." 1600)
